# Shift column C (pulses count) down by one row on sheet "pulseData1":
#   - C1:C200 values move to C2:C201
#   - C1 becomes empty
# This corresponds to selecting C1:C200, cutting, and pasting into C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pulseData1")

$srcValues = $ws.Range("C1:C200").Value()

$dst = $ws.Range("C2:C201")
$dst.Value() = $srcValues

$ws.Range("C1").ClearContents() | Out-Null

$ws.Range("C2:C201").Select() | Out-Null
